$p = $ppt.ActivePresentation

# --- Slide 5: Title "Quality Measures Track" -> "Clinical Reasoning Track" ---
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Clinical Reasoning Track"

# --- Slide 6: body placeholder edits ---
$s6 = $p.Slides.Item(6)
$tr = $s6.Shapes.Item(2).TextFrame.TextRange

# Insert two new level-2 (IndentLevel=2) bullets ("Early bird by 8/16",
# "Deadline by 8/30") right before the "Reach out to track leads" paragraph.
$paraReachOut = $tr.Paragraphs(3, 1)
$paraReachOut.InsertBefore("Early bird by 8/16" + [char]13 + "Deadline by 8/30" + [char]13) | Out-Null

$tr = $s6.Shapes.Item(2).TextFrame.TextRange
$tr.Paragraphs(3, 1).IndentLevel = 2
$tr.Paragraphs(4, 1).IndentLevel = 2

# Split "Attend track kickoffs" into two runs: "Attend " and "track kickoffs"
$tr = $s6.Shapes.Item(2).TextFrame.TextRange
$paraAttend = $tr.Paragraphs(6, 1)
$secondRun = $tr.Characters($paraAttend.Start + 7, $paraAttend.Length - 7)
$secondRun.Text = $secondRun.Text
